$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Clear existing hyperlinks before rewriting data
$ws.Hyperlinks.Delete()

$rows = @(
  @("2026-02-06 18:45:41", "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5460562", "435", "🔥AI,Ai ◆ツール,開発"),
  @("2026-02-06 18:45:41", "AIの改善", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5487324", "310", "🔥AI,Ai"),
  @("2026-02-06 18:45:41", "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5460563", "220", "◆開発,システム開発 ◇管理"),
  @("2026-02-06 18:45:41", "【急募】Next.js × Expoでアプリ開発仲間を大募集!", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486863", "218", "🔥Next.js ◆開発 ◇アプリ"),
  @("2026-02-06 18:45:41", "【急募】WordPressサイト再構築+LINE・予約連携+顧客管理機能構築|テーマ指定あり|", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486856", "93", "◇サイト ○WordPress"),
  @("2026-02-06 18:45:41", "【急募】FlutterでのSNS風アプリ開発をお願いします(Firebase想定)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5487157", "93", "◆開発 ◇アプリ"),
  @("2026-02-06 18:45:41", "【急募】ガイドと旅行者をつなぐマッチングサイト開発", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5475198", "93", "◆開発 ◇サイト"),
  @("2026-02-06 18:45:41", "【業務改善】訪問業務に特化したスケジュール/介入実績管理Webシステム構築", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486583", "85", "◇業務改善"),
  @("2026-02-06 18:45:41", "【Java/講師】新入社員研修のサブ講師募集", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5477144", "78", "★Java"),
  @("2026-02-06 18:45:41", "初心者向けダンススクールの問い合わせフォームを置き換える/拡張するチャットボット開発", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486851", "75", "◆開発"),
  @("2026-02-06 18:45:41", "【募集】PHP + MySQLでのcron用スクリプト作成依頼", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5487013", "50", "◇MySQL ○PHP"),
  @("2026-02-06 18:45:41", "【長期】寝具ブランドのAmazon・楽天市場 運用代行パートナー募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486471", "25", ""),
  @("2026-02-06 18:45:41", "【急募】ECサービスのメール送信障害調査・改善支援", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5487035", "18", ""),
  @("2026-02-06 18:45:41", "【急募】外部CTOを探しています!", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486956", "18", ""),
  @("2026-02-06 18:45:41", "【急募】SSLエラー解決のための専門家を探しています", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486960", "13", ""),
  @("2026-02-06 18:45:41", "【急募】Klaviyoスパム対策とドメイン解決の専門家募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5486673", "13", ""),
  @("2026-02-06 18:45:41", "【急募】BOXファイルをGASでkintoneに自動同期したい", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5487010", "10", "")
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  if ($row[7] -ne "") {
    $ws.Cells.Item($r, 8).Value = $row[7]
  }
  $r = $r + 1
}

# Re-add hyperlinks for F2:F18 in row order
$r = 2
foreach ($row in $rows) {
  $target = $ws.Cells.Item($r, 6)
  $ws.Hyperlinks.Add($target, $row[5])
  $r = $r + 1
}

# Restore the Hyperlink cell style (Hyperlinks.Add creates a duplicate style otherwise)
$ws.Range("F2:F18").Style = "Hyperlink"

# Column B width: 49 -> 52 (account for Excel COM's +5/6 character padding offset)
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664